$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Constant" / "# of LE's" data table in columns J:K ---
$ws.Range("J1").Value = "Constant"
$ws.Range("K1").Value = "# of LE's"

$jValues = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15)
$kValues = @(0, 0, 9, 0, 9, 9, 17, 0, 9, 9, 22, 9, 20, 16, 16)

for ($i = 0; $i -lt $jValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
    $ws.Cells.Item($row, 11).Value = $kValues[$i]
}

# --- Move the scatter chart from its old spot (around F4:M18) down to A11:H25 ---
$co = $ws.ChartObjects(1)
$co.Left = 12.75
$co.Top = 156
$co.Width = 433.0625
$co.Height = 216

# --- Update the cell selection stored in the worksheet view ---
[void]$ws.Range("E1:F1").Select()
